$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so refreshed price strings (which can
# look like plain numbers, e.g. "99.99") are stored as text, matching
# the existing inline-string cells, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.638.80"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.359.87"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "331.06"
$ws.Range("E5").Value = "  +6.96%  "
$ws.Range("D6").Value = "99.99"
$ws.Range("E6").Value = "  -7.65%  "
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "39.69"
$ws.Range("E10").Value = "  -7.29%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "8.44"
$ws.Range("E12").Value = "  -5.45%  "
$ws.Range("D13").Value = "0.996"
$ws.Range("E13").Value = "  -4.76%  "
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "16.27"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "2.724.43"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "2.368.86"
$ws.Range("E17").Value = "  -4.91%  "
$ws.Range("D18").Value = "42.628.56"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "7.81"
$ws.Range("E19").Value = "  +6.95%  "
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "3.77"
$ws.Range("E21").Value = "  +10.40%  "
$ws.Range("D22").Value = "75.46"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "269.74"
$ws.Range("E23").Value = "  +6.84%  "
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -9.21%  "
$ws.Range("D25").Value = "9.97"
$ws.Range("E25").Value = "  +10.87%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "11.45"
$ws.Range("E27").Value = "  -4.34%  "
$ws.Range("D28").Value = "23.20"
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").Value = "175.20"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "0.0902"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").Value = "35.38"
$ws.Range("E33").Value = "  -9.37%  "
$ws.Range("D34").Value = "6.05"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  -8.16%  "
$ws.Range("D37").Value = "0.0359"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("D38").Value = "2.89"
$ws.Range("E38").Value = "  +7.39%  "
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").Value = "3.82"
$ws.Range("E40").Value = "  -6.63%  "
$ws.Range("D41").Value = "1.52"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").Value = "0.233"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "69.82"
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "117.46"
$ws.Range("E45").Value = "  +6.69%  "
$ws.Range("D46").Value = "89.72"
$ws.Range("E46").Value = "  +29.40%  "
$ws.Range("D47").Value = "11.95"
$ws.Range("E47").Value = "  -4.31%  "
$ws.Range("D48").Value = "5.47"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").Value = "9.10"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").Value = "1.570.93"
$ws.Range("E50").Value = "  +5.27%  "
$ws.Range("E51").Value = "  -1.99%  "

# Restore the default (unstyled) cell style now that the text values are set.
$ws.Range("D2:D51").Style = "Normal"
